$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.118.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.744.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5337"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2794"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06181"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.748.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6535"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.645"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.117.51"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006796"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.972.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.435"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.769"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.263"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.516"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  +1.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.780"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("E29").Value = "  +3.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08494"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.803"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.689"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04614"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.658"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9987"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6267"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01621"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.959"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3921"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7488"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.36%  "

$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1150"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.333"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05337"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3502"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.616"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
